# Re-sort the sheet tabs: the sheet named "总计" (currently last) moves to the
# first position, and the sheet named "2022-Q1" (currently first) moves to the
# last position, while "2022-Q2" stays in the middle. This mirrors the
# "resort sheetname" commit: tab order becomes 总计, 2022-Q2, 2022-Q1.

$wb = $excel.ActiveWorkbook

# Move "总计" to be right before the current first sheet ("2022-Q1").
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetQ1 = $wb.Worksheets.Item("2022-Q1")
$sheetTotal.Move($sheetQ1)

# Now move "2022-Q1" to the end, after the last sheet ("2022-Q2").
$sheetQ1 = $wb.Worksheets.Item("2022-Q1")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetQ1.Move($null, $lastSheet)

# "2022-Q1" was the originally-selected/active tab; keep it active after the move.
$sheetQ1 = $wb.Worksheets.Item("2022-Q1")
$sheetQ1.Activate()
